$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

# Writing "01-07-2021" straight into .Value would be auto-parsed by Excel
# into a date serial number (since it looks like a date). To keep it as a
# literal text label (matching the other "Serie" column entries), compute
# it via a text-literal formula on a scratch cell, then copy/paste the
# resulting value (not the formula) into the destination cell. This avoids
# touching NumberFormat/Style on the destination (or anywhere else that
# would stick around), so no new cell style gets introduced.
$scratch = $ws.Cells.Item(1, 20)
$scratch.Formula = '="01-07-2021"'
$scratch.Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item($row, 2).Value = 74463
$ws.Cells.Item($row, 3).Value = 13234
$ws.Cells.Item($row, 4).Value = 7115
$ws.Cells.Item($row, 5).Value = 6565
$ws.Cells.Item($row, 6).Value = 4301
$ws.Cells.Item($row, 7).Value = 6019
$ws.Cells.Item($row, 8).Value = 10177
$ws.Cells.Item($row, 9).Value = 15809
$ws.Cells.Item($row, 10).Value = 11243
